$d = $word.ActiveDocument

$pairs = @(
    @("2024-07-28 Sunday", "2024-07-29 Monday"),
    @("93÷2=", "60÷4="),
    @("25÷4=", "72÷9="),
    @("97÷7=", "64÷2="),
    @("53÷8=", "73÷6="),
    @("22÷6=", "27÷7="),
    @("56÷2=", "93÷6="),
    @("86÷9=", "51÷5="),
    @("26÷3=", "81÷5="),
    @("61÷4=", "72÷5="),
    @("92÷6=", "83÷4="),
    @("23÷7=", "58÷2="),
    @("99÷5=", "89÷2="),
    @("51÷4=", "47÷8="),
    @("44÷9=", "94÷2="),
    @("82÷2=", "26÷5="),
    @("37÷7=", "97÷2="),
    @("68÷9=", "29÷2="),
    @("31÷5=", "32÷8="),
    @("66÷7=", "48÷5="),
    @("69÷2=", "54÷2="),
    @("36÷4=", "53÷6="),
    @("17÷4=", "86÷2="),
    @("44÷2=", "38÷3="),
    @("97÷3=", "87÷8="),
    @("17÷5=", "36÷6=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
